$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for monetization flag
$ws.Range("H1").Value = "monetization"

# Row 2: updated / new data values
$ws.Range("A2").Value = "cringe"
$ws.Range("B2").Value = "E:\storage2\901.mp4"
$ws.Range("C2").Value = "ádasdasdasd"
$ws.Range("G2").Value = "C:\Users\Admin\Downloads\901.mp4"
# force text (not boolean) so it mirrors the source inlineStr "False"
$ws.Range("H2").Value = "'False"

# I2 no longer exists after the edit - clear it so the sheet's used
# range (dimension) shrinks back from A1:I2 to A1:H2
$ws.Range("I2").ClearContents()
